# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.018.72"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "3.766.97"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'602.00"
$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "'165.90"
$ws.Range("E6").Value = "  -2.16%  "

$ws.Range("D7").Value = "3.764.41"
$ws.Range("E7").Value = "  +0.45%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("E10").Value = "  +5.21%  "

$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").Value = "'0.459"
$ws.Range("E12").Value = "  -0.60%  "

$ws.Range("D13").Value = "'37.79"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").Value = "'0.0000248"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").Value = "4.400.46"
$ws.Range("E15").Value = "  +0.69%  "

$ws.Range("D16").Value = "3.754.41"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "69.212.79"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").Value = "'17.75"
$ws.Range("E19").Value = "  +3.35%  "

$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").Value = "'11.33"
$ws.Range("E21").Value = "  +5.25%  "

$ws.Range("D22").Value = "'490.99"
$ws.Range("E22").Value = "  -0.90%  "

$ws.Range("D23").Value = "'0.727"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").Value = "'84.79"
$ws.Range("E25").Value = "  -0.81%  "

$ws.Range("E26").Value = "  -2.01%  "

$ws.Range("E27").Value = "  -0.49%  "

$ws.Range("D28").Value = "'10.08"
$ws.Range("E28").Value = "  -1.35%  "

$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'2.97"
$ws.Range("E30").Value = "  -0.35%  "

$ws.Range("D31").Value = "'8.16"
$ws.Range("E31").Value = "  +2.88%  "

$ws.Range("D32").Value = "'2.44"
$ws.Range("E32").Value = "  -3.94%  "

$ws.Range("D33").Value = "'31.85"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").Value = "3.908.90"

$ws.Range("D35").Value = "3.706.27"
$ws.Range("E35").Value = "  +0.68%  "

$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").Value = "'5.95"
$ws.Range("E37").Value = "  +1.76%  "

$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("D39").Value = "'0.138"
$ws.Range("E39").Value = "  +4.03%  "

$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.12"
$ws.Range("E41").Value = "  +8.36%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'48.61"
$ws.Range("E44").Value = "  -0.83%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'426.81"
$ws.Range("E45").Value = "  -2.65%  "

$ws.Range("D46").Value = "'8.45"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("D48").Value = "'40.23"
$ws.Range("E48").Value = "  -1.37%  "

$ws.Range("E49").Value = "  +10.04%  "

$ws.Range("D50").Value = "'141.55"
$ws.Range("E50").Value = "  +0.28%  "

$ws.Range("D51").Value = "2.794.33"
$ws.Range("E51").Value = "  -0.73%  "
